# Entrega-Final-Cumplimiento-Consignas.xlsx
# "Revision de funcionalidades de responsive"
#
# The commit only touches the "Comentario" (column C) narrative text of a
# handful of rows -- extending/rewording the evidence notes -- plus the
# wrap/vertical-alignment formatting those longer notes need, a couple of
# row heights that grew to fit the new text, and the sheet's saved
# scroll/selection position. Shared-string table churn in the raw XML is
# just Excel's natural side effect of editing cell text/removing no-longer
# referenced strings; we don't need to hand-manage that table ourselves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# ---------------------------------------------------------------------
# 1) Updated "Comentario" text (column C) for several consignas.
# ---------------------------------------------------------------------

$ws.Range("C5").Value = "Se incluye readme.md. Puede encontrarse en el repositorio GIT"

$ws.Range("C7").Value = "Puede encontrarse en static/contacto.html" + $nl + `
    "Cuando se genera el contacto Formspree devuelve el mensaje de recepcion y envia un email a la cuenta registrada"

$ws.Range("C12").Value = "Para netifly: https://sportsprook.netlify.app/ " + $nl + `
    "Para repositorio github:  https://github.com/contesl/"

$ws.Range("C20").Value = "Se implementa en la pagina static/contacto.html" + $nl + `
    "El soporte a multiples formatos esta implementado en css/styles.css"

$ws.Range("C22").Value = "Se implementa en todas las paginas desde los los estilos existentes en css/style.css"

$ws.Range("C24").Value = "Se implemeta el menu de navegacion de Bootstrap." + $nl + `
    "Revisar cada pagina para ver la implementacion"

$ws.Range("C25").Value = "Para repositorio github:  https://github.com/contesl/"

$ws.Range("C27").Value = "Revisar el directorio js para encontrar los scrips que se utilizan en las paginas que asi lo requieren:" + $nl + `
    "productos.js : para el armado dinamico de la pagina de productos luego de un fetch a una api publicada sumado al tratamiento del carrito de compra" + $nl + `
    "resenials.js: para el armado de las reseñas de productos a partir del fetch de la api utilizada en productos.js que tambien contiene las reseñas" + $nl + `
    "validarContacto.js: validacion de datos del formulario de contacto"

$ws.Range("C33").Value = "Esto puede encontrarse en js/productos.js" + $nl + `
    "En la pagina productos.html, cuando se posiciona del cursor en la imagen de un producto, cambia el formato del cursor y se puede acceder a su descripcion ampliada"

$ws.Range("C36").Value = "Esto puede encontrarse en js/productos.js" + $nl + `
    "esto se puede encontrar buscando el comentario   // Cargar el JSON desde el archivo local"

$ws.Range("C39").Value = "Esto puede encontrarse en js/productos.js y  en js/resenias.js porque se utilizo la misma api publica dado que contienen los datos requeridos "

# ---------------------------------------------------------------------
# 2) Formatting follow-up: cells whose notes now wrap across multiple
#    lines need WrapText, and a few vertical alignments shift between
#    top/center to better frame the longer text blocks.
# ---------------------------------------------------------------------

$xlTop = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop
$xlCenter = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignCenter

$ws.Range("C7").WrapText = $true
$ws.Range("C7").VerticalAlignment = $xlTop

$ws.Range("C20").WrapText = $true
$ws.Range("C20").VerticalAlignment = $xlTop

$ws.Range("C24").WrapText = $true
$ws.Range("C24").VerticalAlignment = $xlCenter

$ws.Range("A27").VerticalAlignment = $xlTop

$ws.Range("C33").WrapText = $true
$ws.Range("C33").VerticalAlignment = $xlCenter

$ws.Range("C36").WrapText = $true
$ws.Range("C36").VerticalAlignment = $xlCenter

$ws.Range("A38").VerticalAlignment = $xlTop

$ws.Range("A39").VerticalAlignment = $xlTop

$ws.Range("C39").WrapText = $true
$ws.Range("C39").VerticalAlignment = $xlTop

# ---------------------------------------------------------------------
# 3) Row heights that grew/shrank to accommodate the revised text.
# ---------------------------------------------------------------------

$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 90
$ws.Rows.Item(33).RowHeight = 45
$ws.Rows.Item(39).RowHeight = 39.75

# ---------------------------------------------------------------------
# 4) Saved scroll position / active selection moved further down the
#    sheet (reviewer had scrolled to the bottom rows).
# ---------------------------------------------------------------------

$window = $ws.Application.ActiveWindow
$window.ScrollRow = 27
$ws.Range("C42").Select()
